$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 35, shifting existing rows 35:135 down to 36:136.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new record.
$ws.Range("A35").Value = 5
$ws.Range("B35").Value = "Macroferia Regional de Talca"
$ws.Range("C35").Value = "Maule"
$ws.Range("D35").Value = 44526
$ws.Range("E35").Value = 7
$ws.Range("F35").Value = 100112021
$ws.Range("G35").Value = "Ají"
$ws.Range("H35").Value = "Americana (o)"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 150
$ws.Range("K35").Value = 18000
$ws.Range("L35").Value = 18000
$ws.Range("M35").Value = 18000
$ws.Range("N35").Value = "$/caja 15 kilos"
$ws.Range("O35").Value = "Región del Maule"
$ws.Range("P35").Value = 1200
$ws.Range("Q35").Value = 15
$ws.Range("R35").Value = "Hortaliza"
